# Re-registration test case: duplicate "TestCase01" into a new sheet
# "TestCase01_1" placed right after it, matching how the fixture was
# re-captured for the re-registration scenario.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCase01")

# Duplicate the sheet (Copy places it right after $ws1) so the new sheet
# starts with the same values, number formats and hyperlinks.
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TestCase01_1"

# Approximate the re-captured column widths on the new sheet (columns
# B:D were resized/auto-fit when the data was re-entered).
$ws2.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 32.666666666666664
$ws2.Columns.Item(4).ColumnWidth = 11.75

# The re-captured sheet used the plain default row height (12.75) instead
# of the taller rows on the original sheet.
for ($r = 1; $r -le 4; $r++) {
    $ws2.Rows.Item($r).RowHeight = 12.75
}

# Original sheet is no longer the selected tab; its selection becomes the
# full data range with no distinct active cell offset.
$ws1.Range("A1:D4").Select() | Out-Null

# New sheet is now the active / selected tab, with the active cell left
# just past the data (D5).
$ws2.Select()
$ws2.Range("D5").Select() | Out-Null

$wb.Save()
